$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D9").Value = "공대는 학부만해도 되는데, 통계학은 대학원을 가야된다?"
$ws.Range("E9").Value = "https://blog.pabii.co.kr/engineers-stats-academic-requirements/#utm_source=rss&utm_medium=rss&utm_campaign=engineers-stats-academic-requirements"

$ws.Range("D26").Value = "ai plus(est soft)"

$ws.Range("D37").Value = "[Paper Review] Improving language models by retrieving from trillions of tokens"
$ws.Range("E37").Value = "http://dsba.korea.ac.kr/seminar/?uid=1976&mod=document&pageid=1"

$ws.Range("D52").Value = "[PYTHON] 파이썬(3.8) 표준 라이브러리의 패키지 목록(알파벳 순)"
$ws.Range("E52").Value = "http://ds.sumeun.org/?p=2591&utm_source=rss&utm_medium=rss&utm_campaign=python-%ed%8c%8c%ec%9d%b4%ec%8d%ac3-8-%ed%91%9c%ec%a4%80-%eb%9d%bc%ec%9d%b4%eb%b8%8c%eb%9f%ac%eb%a6%ac%ec%9d%98-%ed%8c%a8%ed%82%a4%ec%a7%80-%eb%aa%a9%eb%a1%9d%ec%95%8c%ed%8c%8c%eb%b2%b3-%ec%88%9c"
